# RDCC-2470 handle trimming in role mapping file
#
# The "BBA9" service id in row 5 (column A) previously carried stray
# leading/trailing whitespace (" BBA9 "). The fix introduces a trimmed
# value ("BBA9 ") and adds a duplicate data row (row 6) using that
# trimmed value paired with the same Role/IDAM Roles as row 5, so the
# mapping file now handles both the old and the new (trimmed) lookup
# key.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: replace " BBA9 " with the trimmed "BBA9 " value.
$ws.Range("A5").Value = "BBA9 "

# New row 6: duplicate of row 5 (Service ID / Role / IDAM Roles).
$ws.Range("A6").Value = "BBA9 "
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = " caseworker-iac-bulkscan "

# Widen column C to fit the data, dropping the old bestFit auto-size flag.
$ws.Columns.Item(3).ColumnWidth = 21.5

# Move the active selection to the newly added row.
$ws.Range("A6").Select()
